$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The weekly topics in column C got shifted down one slot starting the week
# of Oct 28 (row 10), and the old "Java vs. Python reprise and review for
# final exam" entry got split into two separate week topics. The "Nov 18 -
# Nov 22:" date label also lost its trailing colon.

$ws.Range("C10").Value = "GUI Output"
$ws.Range("C12").Value = "GUI Input, Event-driven programming, and MVC"
$ws.Range("B13").Value = "Nov 18 - Nov 22"
$ws.Range("C13").Value = "Decomposition: top down and bottom up"
$ws.Range("C14").Value = "Lambda functions"
$ws.Range("C15").Value = "Java vs. Python"
$ws.Range("C16").Value = "review for final exam"

# The "wrap text" row that used to sit on row 13 (Nov 18 - Nov 22 week) moves
# down to row 14 (Nov 25 - Nov 27 week) along with the shifted content.
$ws.Range("C13").Style = "Normal"
$ws.Rows.Item(13).AutoFit()

$ws.Range("C14").WrapText = $true
$ws.Rows.Item(14).RowHeight = 17

$ws.Range("C15").Style = "Normal"
$ws.Rows.Item(15).AutoFit()

# The saved selection in the sheet view moves to C14.
$ws.Range("C14").Select()
